# Update Name of Algo
# Apply updated numeric values to result_data_KNN worksheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.336
$ws.Range("A3").Value = -21.391
$ws.Range("B5").Value = 6.078
$ws.Range("E5").Value = 12.653
$ws.Range("E9").Value = 12.799
$ws.Range("E11").Value = 13.132
$ws.Range("A14").Value = -20.891
$ws.Range("A21").Value = -20.993
$ws.Range("E21").Value = 13.357
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.27
